$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Drop the now-unused "Density" column (C) entirely.
$ws.Columns.Item(3).Delete()

# Rename "Population" label to lowercase "population".
$ws.Range("A3").Value = "population"

# Add a new row for density, using the value that used to live in column C.
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 995.2500813798443
